$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (2007年 data row), shifting all subsequent rows up by one.
$ws.Rows.Item(2).Delete()
